# Weekly update: prepend two new "Camote" price records to the data table.
# The existing data (rows 59:83) is pushed down by two rows (to 61:85) and a
# pair of brand-new rows are written in their place at rows 59:60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing rows from 59 downward by inserting two blank rows at 59:60.
$ws.Rows("59:60").Insert()

# New row 59
$ws.Cells.Item(59,1).Value  = 9
$ws.Cells.Item(59,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(59,3).Value  = "Metropolitana"
$ws.Cells.Item(59,4).Value  = 44704
$ws.Cells.Item(59,5).Value  = 13
$ws.Cells.Item(59,6).Value  = 100114002
$ws.Cells.Item(59,7).Value  = "Camote"
$ws.Cells.Item(59,8).Value  = "Sin especificar"
$ws.Cells.Item(59,9).Value  = "Primera"
$ws.Cells.Item(59,10).Value = 650
$ws.Cells.Item(59,11).Value = 10000
$ws.Cells.Item(59,12).Value = 10000
$ws.Cells.Item(59,13).Value = 10000
$ws.Cells.Item(59,14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(59,15).Value = "Perú"
$ws.Cells.Item(59,16).Value = 556
$ws.Cells.Item(59,17).Value = 18
$ws.Cells.Item(59,18).Value = "Hortaliza"

# New row 60
$ws.Cells.Item(60,1).Value  = 9
$ws.Cells.Item(60,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(60,3).Value  = "Metropolitana"
$ws.Cells.Item(60,4).Value  = 44704
$ws.Cells.Item(60,5).Value  = 13
$ws.Cells.Item(60,6).Value  = 100114002
$ws.Cells.Item(60,7).Value  = "Camote"
$ws.Cells.Item(60,8).Value  = "Sin especificar"
$ws.Cells.Item(60,9).Value  = "Primera"
$ws.Cells.Item(60,10).Value = 1100
$ws.Cells.Item(60,11).Value = 7500
$ws.Cells.Item(60,12).Value = 8000
$ws.Cells.Item(60,13).Value = 7773
$ws.Cells.Item(60,14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(60,15).Value = "Perú"
$ws.Cells.Item(60,16).Value = 432
$ws.Cells.Item(60,17).Value = 18
$ws.Cells.Item(60,18).Value = "Hortaliza"
